# Fruta / hortaliza, semanal
# Insert two new weekly-report rows at the top of the Frutilla data block
# (new rows 369-370), pushing the existing rows 369-376 down to 371-378.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows 369:376 down by two rows, keeping all of their
# original formatting/values intact (the rows simply shift to 371:378).
$ws.Rows("369:370").Insert()

# Fill in the two brand-new rows (369 and 370) with this week's data,
# copying the date cell's number format from the row right below it.
$ws.Cells.Item(369, 4).NumberFormat = $ws.Cells.Item(371, 4).NumberFormat
$ws.Cells.Item(370, 4).NumberFormat = $ws.Cells.Item(371, 4).NumberFormat

$ws.Cells.Item(369, 1).Value  = 7
$ws.Cells.Item(369, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(369, 3).Value  = "Ñuble"
$ws.Cells.Item(369, 4).Value  = 44890
$ws.Cells.Item(369, 5).Value  = 16
$ws.Cells.Item(369, 6).Value  = "Fruta"
$ws.Cells.Item(369, 7).Value  = 100101
$ws.Cells.Item(369, 8).Value  = "Berries"
$ws.Cells.Item(369, 9).Value  = 100112025
$ws.Cells.Item(369, 10).Value = "Frutilla"
$ws.Cells.Item(369, 11).Value = "Sin especificar"
$ws.Cells.Item(369, 12).Value = "Primera"
$ws.Cells.Item(369, 13).Value = 160
$ws.Cells.Item(369, 14).Value = 6500
$ws.Cells.Item(369, 15).Value = 7000
$ws.Cells.Item(369, 16).Value = 6750
$ws.Cells.Item(369, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(369, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(369, 19).Value = 964
$ws.Cells.Item(369, 20).Value = 7

$ws.Cells.Item(370, 1).Value  = 7
$ws.Cells.Item(370, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(370, 3).Value  = "Ñuble"
$ws.Cells.Item(370, 4).Value  = 44890
$ws.Cells.Item(370, 5).Value  = 16
$ws.Cells.Item(370, 6).Value  = "Fruta"
$ws.Cells.Item(370, 7).Value  = 100101
$ws.Cells.Item(370, 8).Value  = "Berries"
$ws.Cells.Item(370, 9).Value  = 100112025
$ws.Cells.Item(370, 10).Value = "Frutilla"
$ws.Cells.Item(370, 11).Value = "Sin especificar"
$ws.Cells.Item(370, 12).Value = "Segunda"
$ws.Cells.Item(370, 13).Value = 80
$ws.Cells.Item(370, 14).Value = 5000
$ws.Cells.Item(370, 15).Value = 5000
$ws.Cells.Item(370, 16).Value = 5000
$ws.Cells.Item(370, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(370, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(370, 19).Value = 714
$ws.Cells.Item(370, 20).Value = 7

$ws.Range("A1").Select()
